# ---------------------------------------------------------------------------
# Target change (from the OOXML diff / commit "build: move title to sub
# folder"): five <w:abstractNum> entries in word/numbering.xml get a new
# <w:nsid w:val="..."/> GUID (abstractNumId 990, 991, 99411, 99412, 99410).
# Nothing else in the package changes - no paragraph text, no run/paragraph
# formatting, no list appearance (numFmt/lvlText/indent/start values), no
# numId->abstractNumId mapping, no document structure.
#
# w:nsid is the legacy "numbering definition identifier" left over from the
# binary .doc list-tracking format. Word regenerates it internally whenever
# it mints a list definition, but it is not a property any Word
# automation surface (VBA/COM "List", "ListTemplate", "ListFormat", ...)
# ever reads or writes - real Word does not let a macro set it, and this
# host's object model mirrors that: Lists/ListTemplates/ListFormat expose
# visible formatting (NumberFormat, NumberStyle, StartAt, alignment, ...)
# plus the per-document w:numId (List.ListID), never the abstract list's
# w:nsid. There is no Find/Replace route either, since w:nsid never
# appears in any story's visible text - it only lives as an attribute on
# a numbering-part element that no Range/Selection ever addresses.
#
# So this particular edit has no representation in the Word COM object
# model that is exposed to us here; touching list formatting through that
# model (even "setting" a level's NumberFormat back to its own value)
# does not regenerate or otherwise change w:nsid - it is preserved as-is.
# The faithful, risk-free action through this interface is therefore to
# leave the document's object model untouched, which reproduces every
# part of the package byte-for-byte except for container metadata Word
# itself controls on save.
$d = $word.ActiveDocument
